$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: K3 "coremk" -> "coremk_ch2"
$ws.Range("K3").Value = "coremk_ch2"

# Note update: R11 "UPDATED 20180305" -> "UPDATED 20180327"
$ws.Range("R11").Value = "UPDATED 20180327"

# Data value updates (K/N/O raw value columns; dependent OFFSET-formula
# cells in rows 56-63 and 67-74 recalc automatically from these).
$ws.Range("K5").Value = 119.9
$ws.Range("N5").Value = 105.1
$ws.Range("O5").Value = 129.69999999999999
$ws.Range("K6").Value = 8.5
$ws.Range("N6").Value = 3.9
$ws.Range("K7").Value = 76
$ws.Range("N7").Value = 37.9
$ws.Range("O7").Value = 80.5
$ws.Range("K8").Value = 5.4
$ws.Range("N8").Value = 6.6
$ws.Range("O8").Value = 2.4
$ws.Range("K9").Value = 30
$ws.Range("N9").Value = 56.6
$ws.Range("O9").Value = 38.1
$ws.Range("K11").Value = 99.8
$ws.Range("N11").Value = 92.7
$ws.Range("O11").Value = 113
$ws.Range("K12").Value = 9.6999999999999993
$ws.Range("N12").Value = 3.9
$ws.Range("K13").Value = 52.2
$ws.Range("N13").Value = 29.9
$ws.Range("O13").Value = 64
$ws.Range("K14").Value = 8.3000000000000007
$ws.Range("N14").Value = 5.3
$ws.Range("O14").Value = 3.8
$ws.Range("K15").Value = 29.6
$ws.Range("N15").Value = 53.5
$ws.Range("O15").Value = 37.5
$ws.Range("K17").Value = 88.4
$ws.Range("N17").Value = 79.8
$ws.Range("O17").Value = 98.7
$ws.Range("K18").Value = 9.5
$ws.Range("N18").Value = 3.3
$ws.Range("K19").Value = 39.1
$ws.Range("N19").Value = 19.8
$ws.Range("O19").Value = 50.3
$ws.Range("K20").Value = 10.199999999999999
$ws.Range("N20").Value = 4.8
$ws.Range("K21").Value = 29.7
$ws.Range("N21").Value = 51.9
$ws.Range("O21").Value = 37
$ws.Range("K23").Value = 79.900000000000006
$ws.Range("N23").Value = 77.3
$ws.Range("O23").Value = 88.4
$ws.Range("K24").Value = 9.1999999999999993
$ws.Range("N24").Value = 3.3
$ws.Range("K25").Value = 31.3
$ws.Range("N25").Value = 17.899999999999999
$ws.Range("O25").Value = 40.200000000000003
$ws.Range("N26").Value = 4.5999999999999996
$ws.Range("O26").Value = 4.9000000000000004
$ws.Range("K27").Value = 29.7
$ws.Range("N27").Value = 51.5
$ws.Range("O27").Value = 37
$ws.Range("K29").Value = 70.599999999999994
$ws.Range("N29").Value = 74.099999999999994
$ws.Range("O29").Value = 80.099999999999994
$ws.Range("K30").Value = 6.3
$ws.Range("N30").Value = 2.7
$ws.Range("O30").Value = 5.7
$ws.Range("K31").Value = 25.1
$ws.Range("N31").Value = 16.100000000000001
$ws.Range("O31").Value = 31.6
$ws.Range("K32").Value = 9.6
$ws.Range("N32").Value = 4.3
$ws.Range("K33").Value = 29.5
$ws.Range("N33").Value = 51.1
$ws.Range("O33").Value = 36.9
$ws.Range("K35").Value = 68
$ws.Range("N35").Value = 72.900000000000006
$ws.Range("O35").Value = 74.599999999999994
$ws.Range("K36").Value = 6.3
$ws.Range("N36").Value = 3.3
$ws.Range("K37").Value = 22.4
$ws.Range("N37").Value = 14.6
$ws.Range("O37").Value = 25.6
$ws.Range("K38").Value = 9.8000000000000007
$ws.Range("N38").Value = 4.3
$ws.Range("K39").Value = 29.5
$ws.Range("N39").Value = 50.7
$ws.Range("O39").Value = 36.9
$ws.Range("K41").Value = 60.6
$ws.Range("N41").Value = 71.400000000000006
$ws.Range("O41").Value = 68.7
$ws.Range("K42").Value = 2.5
$ws.Range("N42").Value = 2.6
$ws.Range("O42").Value = 4.9000000000000004
$ws.Range("K43").Value = 18.899999999999999
$ws.Range("N43").Value = 13.9
$ws.Range("K44").Value = 9.6
$ws.Range("N44").Value = 4.3
$ws.Range("O44").Value = 6.8
$ws.Range("K45").Value = 29.5
$ws.Range("N45").Value = 50.5
$ws.Range("O45").Value = 36.9
$ws.Range("K47").Value = 58.9
$ws.Range("N47").Value = 70.2
$ws.Range("O47").Value = 67
$ws.Range("K48").Value = 2.5
$ws.Range("N48").Value = 1.7
$ws.Range("O48").Value = 4.3
$ws.Range("K49").Value = 17.100000000000001
$ws.Range("N49").Value = 13.5
$ws.Range("O49").Value = 18.899999999999999
$ws.Range("N50").Value = 4.4000000000000004
$ws.Range("O50").Value = 6.9
$ws.Range("K51").Value = 29.6
$ws.Range("N51").Value = 50.5
$ws.Range("O51").Value = 36.9

# Fill color used by the R11:S11 note highlight changes from orange to red
$ws.Range("R11:S11").Interior.Color = 255

# Selection / view state
$ws.Range("Q28").Select()
